$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 804.20636
$ws.Range("J132").Value = 515
$ws.Range("L132").Value = 1545
$ws.Range("N132").Value = -6605
$ws.Range("H137").Value = 1572.1428
$ws.Range("I137").Value = 1442.4
$ws.Range("J137").Value = 1644.2222
$ws.Range("K137").Value = 4327.200000000001
$ws.Range("L137").Value = 4932.6666
$ws.Range("M137").Value = -1777.200000000001
$ws.Range("N137").Value = -10032.6666
$ws.Range("H138").Value = 3737.76
$ws.Range("I138").Value = 1887.5
$ws.Range("J138").Value = 3814.8542
$ws.Range("K138").Value = 5662.5
$ws.Range("L138").Value = 11444.5626
$ws.Range("M138").Value = -522.5
$ws.Range("N138").Value = -21724.5626

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1480.7646
$ws.Range("I2").Value = 1032.6923
$ws.Range("J2").Value = 2937
$ws.Range("K2").Value = 1032.6923
$ws.Range("L2").Value = 2937
$ws.Range("M2").Value = -919.6922999999999
$ws.Range("N2").Value = -3163
$ws.Range("H45").Value = 4053
$ws.Range("I45").Value = 3748.5
$ws.Range("K45").Value = 3748.5
$ws.Range("M45").Value = -3371.5
$ws.Range("H61").Value = 2834
$ws.Range("I61").Value = 2567.5667
$ws.Range("J61").Value = 4166.1665
$ws.Range("K61").Value = 2567.5667
$ws.Range("L61").Value = 4166.1665
$ws.Range("M61").Value = -2355.5667
$ws.Range("N61").Value = -4590.1665
$ws.Range("H74").Value = 2583.3333
$ws.Range("I74").Value = 2617.647
$ws.Range("K74").Value = 2617.647
$ws.Range("M74").Value = -1743.647
$ws.Range("H77").Value = 2583.3333
$ws.Range("I77").Value = 2617.647
$ws.Range("K77").Value = 13088.235
$ws.Range("M77").Value = -8720.235000000001
$ws.Range("H102").Value = 1615.1666
$ws.Range("I102").Value = 1318.5264
$ws.Range("K102").Value = 1318.5264
$ws.Range("M102").Value = 303.4736
$ws.Range("H116").Value = 1480.7646
$ws.Range("I116").Value = 1032.6923
$ws.Range("J116").Value = 2937
$ws.Range("K116").Value = 1032.6923
$ws.Range("L116").Value = 2937
$ws.Range("M116").Value = 1261.3077
$ws.Range("N116").Value = -7525
$ws.Range("H122").Value = 3132.1724
$ws.Range("I122").Value = 2340.0527
$ws.Range("J122").Value = 4637.2
$ws.Range("K122").Value = 7020.158100000001
$ws.Range("L122").Value = 13911.6
$ws.Range("M122").Value = -4570.158100000001
$ws.Range("N122").Value = -18811.6
$ws.Range("H132").Value = 2328.3872
$ws.Range("I132").Value = 2144.138
$ws.Range("K132").Value = 6432.414
$ws.Range("M132").Value = -3902.414
$ws.Range("H136").Value = 2834
$ws.Range("I136").Value = 2567.5667
$ws.Range("J136").Value = 4166.1665
$ws.Range("K136").Value = 7702.7001
$ws.Range("L136").Value = 12498.4995
$ws.Range("M136").Value = -5152.7001
$ws.Range("N136").Value = -17598.4995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1480.7646
$ws.Range("I3").Value = 1032.6923
$ws.Range("J3").Value = 2937
$ws.Range("K3").Value = 1032.6923
$ws.Range("L3").Value = 2937
$ws.Range("M3").Value = -918.6922999999999
$ws.Range("N3").Value = -3165
$ws.Range("H94").Value = 825.1724
$ws.Range("I94").Value = 764.1905
$ws.Range("J94").Value = 985.25
$ws.Range("K94").Value = 764.1905
$ws.Range("L94").Value = 985.25
$ws.Range("M94").Value = -313.1905
$ws.Range("N94").Value = -1887.25
$ws.Range("H105").Value = 1671.7931
$ws.Range("I105").Value = 1311.4546
$ws.Range("K105").Value = 1311.4546
$ws.Range("M105").Value = 435.5454

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1859.2046
$ws.Range("J31").Value = 1864.6757
$ws.Range("L31").Value = 1864.6757
$ws.Range("N31").Value = -2454.6757
$ws.Range("H34").Value = 1859.2046
$ws.Range("J34").Value = 1864.6757
$ws.Range("L34").Value = 1864.6757
$ws.Range("N34").Value = -2268.6757
$ws.Range("H58").Value = 4256.3887
$ws.Range("I58").Value = 3106.4
$ws.Range("J58").Value = 4698.6924
$ws.Range("K58").Value = 3106.4
$ws.Range("L58").Value = 4698.6924
$ws.Range("M58").Value = -2903.4
$ws.Range("N58").Value = -5104.6924
$ws.Range("H105").Value = 1079.7
$ws.Range("I105").Value = 810.7778
$ws.Range("K105").Value = 810.7778
$ws.Range("M105").Value = 936.2222
$ws.Range("H122").Value = 5492.125
$ws.Range("I122").Value = 6179.3335
$ws.Range("J122").Value = 5079.8
$ws.Range("K122").Value = 18538.0005
$ws.Range("L122").Value = 15239.4
$ws.Range("M122").Value = -16088.0005
$ws.Range("N122").Value = -20139.4
$ws.Range("H132").Value = 3653.4443
$ws.Range("I132").Value = 3054.7144
$ws.Range("J132").Value = 5749
$ws.Range("K132").Value = 9164.143199999999
$ws.Range("L132").Value = 17247
$ws.Range("M132").Value = -6634.143199999999
$ws.Range("N132").Value = -22307
$ws.Range("H134").Value = 1571.975
$ws.Range("I134").Value = 1456.2433
$ws.Range("K134").Value = 4368.7299
$ws.Range("M134").Value = -1833.7299
$ws.Range("H136").Value = 4256.3887
$ws.Range("I136").Value = 3106.4
$ws.Range("J136").Value = 4698.6924
$ws.Range("K136").Value = 9319.200000000001
$ws.Range("L136").Value = 14096.0772
$ws.Range("M136").Value = -6769.200000000001
$ws.Range("N136").Value = -19196.0772

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2490.8096
$ws.Range("I132").Value = 2291.8572
$ws.Range("J132").Value = 2590.2856
$ws.Range("K132").Value = 20626.7148
$ws.Range("L132").Value = 23312.5704
$ws.Range("M132").Value = -18096.7148
$ws.Range("N132").Value = -28372.5704
$ws.Range("H137").Value = 2408.889
$ws.Range("I137").Value = 1960
$ws.Range("K137").Value = 5880
$ws.Range("M137").Value = -780

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18388.2
$ws.Range("I70").Value = 44911
$ws.Range("K70").Value = 44911
$ws.Range("M70").Value = -44641
$ws.Range("H73").Value = 18388.2
$ws.Range("I73").Value = 44911
$ws.Range("K73").Value = 44911
$ws.Range("M73").Value = -43975
$ws.Range("H102").Value = 1263.909
$ws.Range("I102").Value = 1216.4814
$ws.Range("J102").Value = 1477.3334
$ws.Range("K102").Value = 1216.4814
$ws.Range("L102").Value = 1477.3334
$ws.Range("M102").Value = 405.5186000000001
$ws.Range("N102").Value = -4721.3334
$ws.Range("H107").Value = 150
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 150
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 150
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -3990
$ws.Range("H113").Value = 19633.723
$ws.Range("I113").Value = 2356
$ws.Range("J113").Value = 46784.43
$ws.Range("K113").Value = 2356
$ws.Range("L113").Value = 46784.43
$ws.Range("M113").Value = -186
$ws.Range("N113").Value = -51124.43
$ws.Range("H122").Value = 3208.889
$ws.Range("I122").Value = 3577.8462
$ws.Range("K122").Value = 10733.5386
$ws.Range("M122").Value = -8283.5386
$ws.Range("H132").Value = 3738.1738
$ws.Range("J132").Value = 4562.5
$ws.Range("L132").Value = 13687.5
$ws.Range("N132").Value = -18747.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2374.875
$ws.Range("I100").Value = 2374.875
$ws.Range("K100").Value = 2374.875
$ws.Range("M100").Value = -1833.875
$ws.Range("H122").Value = 4316.6665
$ws.Range("I122").Value = 4379
$ws.Range("J122").Value = 4005
$ws.Range("K122").Value = 13137
$ws.Range("L122").Value = 12015
$ws.Range("M122").Value = -10687
$ws.Range("N122").Value = -16915
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5565.5
$ws.Range("I62").Value = 5633.5713
$ws.Range("J62").Value = 5497.4287
$ws.Range("K62").Value = 5633.5713
$ws.Range("L62").Value = 5497.4287
$ws.Range("M62").Value = -5009.5713
$ws.Range("N62").Value = -6745.4287
$ws.Range("H65").Value = 5565.5
$ws.Range("I65").Value = 5633.5713
$ws.Range("J65").Value = 5497.4287
$ws.Range("K65").Value = 28167.8565
$ws.Range("L65").Value = 27487.1435
$ws.Range("M65").Value = -25047.8565
$ws.Range("N65").Value = -33727.14350000001
$ws.Range("H122").Value = 76929450
$ws.Range("I122").Value = 90915816
$ws.Range("J122").Value = 4400
$ws.Range("K122").Value = 272747448
$ws.Range("L122").Value = 13200
$ws.Range("M122").Value = -272744998
$ws.Range("N122").Value = -18100
$ws.Range("H132").Value = 1689.6482
$ws.Range("I132").Value = 1541.8864
$ws.Range("J132").Value = 2339.8
$ws.Range("K132").Value = 4625.6592
$ws.Range("L132").Value = 7019.400000000001
$ws.Range("M132").Value = -2095.6592
$ws.Range("N132").Value = -12079.4
